$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) — F column updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 1381
$ws1.Range("F5").Value = 108
$ws1.Range("F7").Value = 11777
$ws1.Range("F8").Value = 4408
$ws1.Range("F9").Value = 31
$ws1.Range("F11").Value = 27
$ws1.Range("F13").Value = 2555
$ws1.Range("F15").Value = 155
$ws1.Range("F17").Value = 5125
$ws1.Range("F20").Value = 524
$ws1.Range("F21").Value = 11361
$ws1.Range("F22").Value = 11308

# Sheet "全部类型" (sheetId 4) — F column updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1381
$ws4.Range("F5").Value = 108
$ws4.Range("F7").Value = 11777
$ws4.Range("F8").Value = 4408
$ws4.Range("F9").Value = 31
$ws4.Range("F11").Value = 27
$ws4.Range("F13").Value = 2555
$ws4.Range("F16").Value = 155
$ws4.Range("F18").Value = 5125
$ws4.Range("F21").Value = 524
$ws4.Range("F22").Value = 11361
$ws4.Range("F23").Value = 11308
